$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete worksheet row 14 entirely; this shifts rows 15-24 up to 14-23
$ws.Rows.Item(14).Delete()
